# small change in load state of ins_cache
#
# The ins_cache rows (DDR locations 0x1010..0x101F, spreadsheet rows 18-33)
# stopped tracking the running "value" column (C) for their ADD/SUB/TSC/ABS/
# FINAL columns (D:I) and instead stay pinned at the cache's last loaded
# value, 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Pin columns D:I (ADD, ADD, SUB, TSC, ABS, FINAL) to 16 for rows 18 through 33.
$ws.Range("D18:I33").Value = 16

# Reflect the author's updated scroll position / selection in the sheet view.
$ws.Activate() | Out-Null
$ws.Range("I50").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 88
$excel.ActiveWindow.ScrollColumn = 1
